# Add a new "ema" column (G) to the BTC price sheet.
# G1 gets the same bold/bordered header style as the other headers (A1:F1),
# G2 seeds the EMA with the first Close value (column E), and every
# subsequent row applies a standard exponential-moving-average update
# with a 7-period smoothing factor (alpha = 2/(7+1) = 0.25):
#   ema[row] = Close[row] * alpha + ema[row-1] * (1 - alpha)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Header cell: clone the formatting already used by the other header cells
# (bold font + thin border, centered/top aligned) so no new style is created.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null
$ws.Range("G1").Value = "ema"

$alpha = 2.0 / (7.0 + 1.0)

$prevEma = $ws.Cells.Item(2, 5).Value2
$ws.Cells.Item(2, 7).Value = $prevEma

for ($row = 3; $row -le $lastRow; $row++) {
    $close = $ws.Cells.Item($row, 5).Value2
    $ema = $close * $alpha + $prevEma * (1 - $alpha)
    $ws.Cells.Item($row, 7).Value = $ema
    $prevEma = $ema
}
